$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-purpose the sheet from the old "URL / Expected Behavior" table to the
# new single-column "Lipid Test" scenario list.
# Set values in this order so the shared-string table indices line up with
# the target workbook (HDL=0, Cholesterol=1, Lipid Test=2, TestName=3).
$ws.Range("A4").Value = "HDL"
$ws.Range("A2").Value = "Cholesterol"
$ws.Range("A3").Value = "Lipid Test"
$ws.Range("A1").Value = "TestName"

# Drop column B entirely - the new layout only uses column A.
$ws.Range("B1:B4").Clear()

# Wrap text for the header + data rows that keep the shared "s=1" style.
$ws.Range("A1:A3").WrapText = $true

# The last row (A4) goes back to the default/no style.
$ws.Range("A4").Style = "Normal"

# Make the header row a bit taller.
$ws.Rows(1).RowHeight = 19.5

# Move the active selection to A5, just below the data.
$ws.Range("A5").Select() | Out-Null
